# "updated activity till excel form"
# Rows 3-15 (runs/balls/fours/sixes in columns C:F) get reshuffled to a new
# set of per-innings figures; row 2 is untouched. Values are written with a
# leading apostrophe so they land as text (matching the rest of the sheet,
# which stores these "numeric" stats as text / numberStoredAsText).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStats = @{
    3  = @("21", "24", "2", "0")
    4  = @("36", "35", "3", "1")
    5  = @("2",  "3",  "0", "0")
    6  = @("3",  "9",  "0", "0")
    7  = @("27", "33", "2", "1")
    8  = @("14", "20", "1", "0")
    9  = @("56", "38", "4", "2")
    10 = @("28", "27", "1", "2")
    11 = @("38", "17", "5", "1")
    12 = @("37", "25", "5", "0")
    13 = @("31", "29", "4", "0")
    14 = @("37", "25", "3", "2")
    15 = @("5",  "9",  "0", "0")
}

foreach ($row in $newStats.Keys) {
    $vals = $newStats[$row]
    $ws.Range("C$row").Value = "'" + $vals[0]
    $ws.Range("D$row").Value = "'" + $vals[1]
    $ws.Range("E$row").Value = "'" + $vals[2]
    $ws.Range("F$row").Value = "'" + $vals[3]
}
